$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper colors (BGR-encoded Long values for Excel COM .Color) ---
# Interior: FFC7CE (light red)  -> R=255 G=199 B=206
$badFill = 13551615
# Font: 9C0006 (dark red)       -> R=156 G=0   B=6
$badFont = 393372

function Set-Bad($addr) {
    $r = $ws.Range($addr)
    $r.ClearContents()
    $r.Interior.Color = $badFill
    $r.Font.Color = $badFont
}

# --- Grid value updates (precedence table body, rows 3-11) ---

# Row 3 ("#" row): G3 switches from ">" to "<"; B3 becomes an empty "Bad" cell
Set-Bad "B3"
$ws.Range("G3").Value = "<"

# Row 6 / Row 7 ("..", "< <= > >= == ~=" rows): F6/F7 "Err" -> ">"
$ws.Range("F6").Value = ">"
$ws.Range("F7").Value = ">"

# Row 8 ("(" row): J8 "Err" -> empty Bad cell; also fix its border (loses
# top/bottom, keeps thin left/right) now that the neighboring Err cell is gone
Set-Bad "J8"
$ws.Range("H8").Borders.Item(8).LineStyle = -4142
$ws.Range("H8").Borders.Item(9).LineStyle = -4142

# Row 9 (")" row): B9 "<" -> empty Bad; G9/I9 "Err" -> empty Bad
Set-Bad "B9"
Set-Bad "G9"
Set-Bad "I9"

# Row 10 ("id" row): B10 ">" -> "<"; G10/I10 "Err" -> empty Bad
$ws.Range("B10").Value = "<"
Set-Bad "G10"
Set-Bad "I10"

# Row 11 ("$" row): H11 "Err" -> empty Bad; J11 "Finish" -> "END"
Set-Bad "H11"
$ws.Range("J11").Value = "END"

# --- Selection state ---
$ws.Range("I21").Select()
